$wb = $excel.ActiveWorkbook

# --- Update cell selections on the affected sheets (use original sheet
# names since the rename happens later in this script) -----------------
$wb.Worksheets.Item("SetPassword").Range("B4").Select() | Out-Null
$wb.Worksheets.Item("CancelVote").Range("F40").Select() | Out-Null
$wb.Worksheets.Item("ViewDuplicateVote").Range("H24").Select() | Out-Null
$wb.Worksheets.Item("ViewProcessingVote").Range("D15").Select() | Out-Null

# Make ViewVoteParticipator the final active sheet/tab with its new
# selection (this also flips the tabSelected flag from SignUp to this
# sheet, and updates the workbook's activeTab index).
$wb.Worksheets.Item("ViewVoteParticipator").Range("I20").Select() | Out-Null

# --- Rename sheets -------------------------------------------------------
$wb.Worksheets.Item("IndexVoteHistory").Name = "IndexVoteHistoryList"
$wb.Worksheets.Item("CancelVote").Name = "CancelUserSelection"
$wb.Worksheets.Item("ViewDuplicateVote").Name = "ViewDuplicateSelection"
